# edit.ps1 - applies the "typo fixed in resume" change:
#   1. Slide 1: "Mobile" -> "Data Science" skill label, with an explicit
#      position/size override (previously inherited from the layout).
#   2. Slide 3: reposition the "Data Science" project-timeframe textbox.
#   3. Slide 4: reposition the "Bot" project-timeframe textbox.

$p = $ppt.ActivePresentation

# --- Slide 1: "Mobile" skill placeholder (id=323) -----------------------
$slide1 = $p.Slides.Item(1)
$mobileShape = $null
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.Id -eq 323) {
        $mobileShape = $shp
    }
}

$mobileShape.Left = 284.2558267716535
$mobileShape.Top = 667.0453797307086
$mobileShape.Width = 55.85755905511811
$mobileShape.Height = 20.745905911811025
$mobileShape.TextFrame.TextRange.Text = "Data Science"

# --- Slide 3: "Data Science" textbox (id=4) reposition -------------------
$slide3 = $p.Slides.Item(3)
$dataScienceBox = $null
for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
    $shp = $slide3.Shapes.Item($i)
    if ($shp.Id -eq 4) {
        $dataScienceBox = $shp
    }
}

$dataScienceBox.Left = 42.7544094488189
$dataScienceBox.Top = 363.2803937007874

# --- Slide 4: "Bot" textbox (id=6) reposition -----------------------------
$slide4 = $p.Slides.Item(4)
$botBox = $null
for ($i = 1; $i -le $slide4.Shapes.Count; $i++) {
    $shp = $slide4.Shapes.Item($i)
    if ($shp.Id -eq 6) {
        $botBox = $shp
    }
}

$botBox.Left = 41.34527559055118
$botBox.Top = 20.62472440944882
